# Updates cryptos list values (Price / Volume(1h) columns, and the
# WEMIXToken / BabyDogeCoin row swap) to match the latest scraped data.
#
# NOTE: Many "Price" values look numeric (e.g. "27.197.54", "0.808") but
# must be stored as literal TEXT, exactly as in the source data (the
# workbook stores every Coin/Link/Price/Volume cell as a string). Setting
# Range.Value directly with such strings causes Excel to auto-convert them
# to numbers/dates, corrupting the data. To avoid this we instead write a
# formula that evaluates to the desired literal text (wrapping any
# non-ASCII character, such as the subscript digits used for very small
# prices, in UNICHAR()), then convert that formula to a plain value via
# Copy + PasteSpecial(xlPasteValues). This guarantees the cell ends up
# holding exactly the intended text, regardless of how it looks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Range, $FormulaLiteral)
    $Range.Formula = "=" + $FormulaLiteral
    $Range.Copy() | Out-Null
    $Range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

Set-CellText $ws.Range("D2") '"27.197.54"'
Set-CellText $ws.Range("E2") '"  +0.16%  "'
Set-CellText $ws.Range("D3") '"1.632.19"'
Set-CellText $ws.Range("E3") '"  -0.73%  "'
Set-CellText $ws.Range("E4") '"  -0.02%  "'
Set-CellText $ws.Range("D5") '"216.72"'
Set-CellText $ws.Range("E5") '"  -0.31%  "'
Set-CellText $ws.Range("E6") '"  +1.46%  "'
Set-CellText $ws.Range("E7") '"  -0.04%  "'
Set-CellText $ws.Range("E8") '"  -0.29%  "'
Set-CellText $ws.Range("E9") '"  -0.85%  "'
Set-CellText $ws.Range("E10") '"  +1.71%  "'
Set-CellText $ws.Range("E11") '"  +0.50%  "'
Set-CellText $ws.Range("D12") '"1.639.36"'
Set-CellText $ws.Range("E12") '"  -0.44%  "'
Set-CellText $ws.Range("E13") '"  +0.12%  "'
Set-CellText $ws.Range("E14") '"  +0.80%  "'
Set-CellText $ws.Range("D15") '"27.177.47"'
Set-CellText $ws.Range("E15") '"  +0.13%  "'
Set-CellText $ws.Range("D16") '"64.74"'
Set-CellText $ws.Range("E16") '"  -4.01%  "'
Set-CellText $ws.Range("D18") '"215.47"'
Set-CellText $ws.Range("E18") '"  -1.38%  "'
Set-CellText $ws.Range("E19") '"  +0.06%  "'
Set-CellText $ws.Range("E20") '"  +1.02%  "'
Set-CellText $ws.Range("E21") '"  -0.76%  "'
Set-CellText $ws.Range("D22") '"2.50"'
Set-CellText $ws.Range("E22") '"  -0.10%  "'
Set-CellText $ws.Range("D23") '"9.09"'
Set-CellText $ws.Range("E23") '"  -0.92%  "'
Set-CellText $ws.Range("D24") '"148.25"'
Set-CellText $ws.Range("E24") '"  +0.48%  "'
Set-CellText $ws.Range("E25") '"  +0.01%  "'
Set-CellText $ws.Range("E26") '"  -1.63%  "'
Set-CellText $ws.Range("E27") '"  -0.05%  "'
Set-CellText $ws.Range("D28") '"15.57"'
Set-CellText $ws.Range("E28") '"  -1.13%  "'
Set-CellText $ws.Range("D29") '"0.0504"'
Set-CellText $ws.Range("E29") '"  -0.11%  "'
Set-CellText $ws.Range("E30") '"  -0.52%  "'
Set-CellText $ws.Range("E31") '"  +0.80%  "'
Set-CellText $ws.Range("E32") '"  -0.68%  "'
Set-CellText $ws.Range("D33") '"1.311.91"'
Set-CellText $ws.Range("E33") '"  +3.55%  "'
Set-CellText $ws.Range("D34") '"1.57"'
Set-CellText $ws.Range("E34") '"  -1.33%  "'
Set-CellText $ws.Range("E35") '"  +0.01%  "'
Set-CellText $ws.Range("E36") '"  -1.36%  "'
Set-CellText $ws.Range("D37") '"0.851"'
Set-CellText $ws.Range("E37") '"  +1.25%  "'
Set-CellText $ws.Range("D38") '"0.541"'
Set-CellText $ws.Range("E38") '"  -0.29%  "'
Set-CellText $ws.Range("E39") '"  -0.04%  "'
Set-CellText $ws.Range("D40") '"2.25"'
Set-CellText $ws.Range("E40") '"  +1.31%  "'
Set-CellText $ws.Range("D41") '"0.805"'
Set-CellText $ws.Range("E41") '"  -0.28%  "'
Set-CellText $ws.Range("D42") '"63.92"'
Set-CellText $ws.Range("E42") '"  +2.68%  "'
Set-CellText $ws.Range("D43") '"1.770.44"'
Set-CellText $ws.Range("E43") '"  -0.86%  "'
Set-CellText $ws.Range("D44") '"5.23"'
Set-CellText $ws.Range("E44") '"  -3.04%  "'
Set-CellText $ws.Range("D45") '"90.74"'
Set-CellText $ws.Range("E45") '"  -1.06%  "'
Set-CellText $ws.Range("E46") '"  -0.30%  "'
# row 47 becomes BabyDogeCoin (was WEMIXToken)
Set-CellText $ws.Range("B47") '"BabyDogeCoin"'
Set-CellText $ws.Range("C47") '"https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"'
Set-CellText $ws.Range("D47") '"0.0"&UNICHAR(8326)&"0106"'
Set-CellText $ws.Range("E47") '"  -1.08%  "'

# row 48 becomes WEMIXToken (was BabyDogeCoin)
Set-CellText $ws.Range("B48") '"WEMIXToken"'
Set-CellText $ws.Range("C48") '"https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"'
Set-CellText $ws.Range("D48") '"0.808"'
Set-CellText $ws.Range("E48") '"  +20.48%  "'

Set-CellText $ws.Range("D49") '"0.0516"'
Set-CellText $ws.Range("E49") '"  +0.98%  "'
Set-CellText $ws.Range("D50") '"7.56"'
Set-CellText $ws.Range("E50") '"  -1.67%  "'
Set-CellText $ws.Range("E51") '"  -0.30%  "'
$excel.CutCopyMode = 0
